# Fixed update to excel issue
# - Rename "Requested quantity" headers to metric-specific names on the
#   existing "Weekly Quantity" / "Monthly Trend" sheets.
# - Add a new "PO Forecast" sheet with forecasted PO quantity data.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
  @(44955.99999999999, 112, 111.9996908626221,  111.9996908626558),
  @(44983.99999999999, 8,   7.999738863077891,  7.999738863111081),
  @(44990.99999999999, 0,  -18.00024916242282,  -18.00024910169242),
  @(44997.99999999999, 0,  -44.00023734871169,  -44.00023693798391),
  @(45004.99999999999, 0,  -70.00022566482752,  -70.000224686165),
  @(45011.99999999999, 0,  -96.00021410476626,  -96.00021235418329),
  @(45018.99999999999, 0,  -122.0002024800792,  -122.0001999675792),
  @(45025.99999999999, 0,  -148.0001909158929,  -148.0001875436351),
  @(45032.99999999999, 0,  -174.0001793872482,  -174.0001750153421),
  @(45039.99999999999, 0,  -200.0001679566864,  -200.0001624202703)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Date formatting on column A to match the other sheets' "ds" columns.
$wsForecast.Range("A2:A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
